# "park pick up money event"
# Adds a new EventType(EV) row pair for a Park work/pickup event:
#   Park-work / EV015 / Action
#   Park      / EV016 / Fade Out Persist
# Renames the CareQuiz EventType(ET) row to "Facility-CareQuiz".

$wb = $excel.ActiveWorkbook

# ---- Sheet: EventType(EV) ----
$ws1 = $wb.Worksheets.Item("EventType(EV)")

# Insert two new rows after the existing "Park/EV014/Text" row (row 15),
# pushing the "Every/EV998" and "OnlyScript/EV999" rows down.
$ws1.Range("A16:C17").Insert()

$ws1.Range("A16").Value = "Park-work"
$ws1.Range("B16").Value = "EV015"
$ws1.Range("C16").Value = "Action"

$ws1.Range("A17").Value = "Park"
$ws1.Range("B17").Value = "EV016"
$ws1.Range("C17").Value = "Fade Out Persist"

# Column A got a touch wider to fit "Park-work".
$ws1.Columns.Item(1).ColumnWidth = 12

# Selection / view bookkeeping to match the authored workbook state.
$ws1.Application.ActiveWindow.ScrollRow = 4
$ws1.Range("C17").Select()

# ---- Sheet: EventType(ET) ----
$ws2 = $wb.Worksheets.Item("EventType(ET)")

# The CareQuiz event type becomes facility-specific.
$ws2.Range("A2").Value = "Facility-CareQuiz"

$ws2.Columns.Item(1).ColumnWidth = 15.09765625
$ws2.Range("B3").Select()

# ---- Sheet: ChoiceEvent ----
$ws3 = $wb.Worksheets.Item("ChoiceEvent")
$ws3.Range("D8").Select()
